# Sprint 1 Burndown Chart - update actuals ("Storypoints Ist") data.
# Fills in the daily "erledigt" (done) storypoints for column B (Soll/Ist
# verbleibend calc source) and column D (Storypoints Plan erledigt), which
# were previously blank for rows 4-6 (except B4 which had a placeholder of
# 7 that becomes 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
$ws.Range("B6").Value = 3

$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 3

# Match the author's final selection (cell B4) as recorded in the sheet view.
$ws.Range("B4").Select() | Out-Null
